$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 24 to the "Logs" sheet ---
$logs.Range("A24").Value = "Zorg jij dat deze pallets vandaag nog klaarstaan?"
$logs.Range("B24").Value = "mailmind.test@zohomail.eu"
$logs.Range("C24").Value = "Testmail #14: Zorg jij dat deze pallets vandaag nog klaarstaan?"
$logs.Range("D24").Value = "Planning / Afspraak"
$logs.Range("E24").Value = "Geachte afzender,
Hartelijk dank voor uw e-mail. Helaas is deze e-mail niet voor ons bedoeld, aangezien wij een professionele e-mailassistent zijn en geen fysieke pallets kunnen klaarzetten. Voor verdere assistentie met betrekking tot pallets verzoek ik u vriendelijk om contact op te nemen met de desbetreffende afdeling binnen uw organisatie.
Met vriendelijke groet,
[Bedrijfsnaam] E-mailassistent"
$logs.Range("F24").Value = "2025-07-23 22:45:05"
$logs.Range("G24").Value = "Ja"
$logs.Range("H24").Value = "Nee"
$logs.Range("I24").Value = "Ja"
$logs.Range("J24").Value = "Nee"
$logs.Rows.Item(24).AutoFit()

# --- Extend the conditional-formatting ranges from row 23 to row 24 ---
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "23")
    $newRange = $logs.Range($col + "2:" + $col + "24")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard summary count for "Planning / Afspraak" ---
$dashboard.Range("B5").Value = 3
